$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headline = 'promising biotech stock to consider july th best travel stock worth watching july th top small cap stock to follow now july th bailard inc buy share of axon enterprise inc nasdaqaxon airline stock to follow now july th how trump tax law will reshape benefit qa is this the end of the magnificent seven boss stop trying to flatter your worker say leadership expertit doesnt make them like you verkkokauppacom oyj manager transaction nina anttila what separate successful manager from the pack come down to their oneonones here how to do it right today quordle hint and answer for july verkkokauppacom oyj manager transaction pekka litmanen sp signal the u stock rally momentum is waning bloomberg what would you like to leave behind a financial planner guide to family wealth discussion prediction stock that will be worth more than annaly capital year from now top budget gaming phone built for speed top pick under r autohome inc to announce second quarter and interim financial result on july massachusetts financial service co ma sell share of quidelortho corporation nasdaqqdel sbi security co ltd purchase share of fabrinet nysefn apollon wealth management llc acquires share of franklin bsp realty trust inc nysefbrt massachusetts financial service co ma ha million position in royal gold inc nasdaqrgld bitcoin miner iren could surge a sotp valuation peg target hedge fund manager eric jackson say were cooking with gas stock to watch monday verizon bp stellantis block she got the last seat on a flight next to a stranger theyve been married for year zwsoft launch zwcad flex in the u and canada powerful cad performance meet purchasing flexibility cwm llc grows position in lamar advertising company nasdaqlamr top chinese stock to add to your watchlist july th bailard inc sell share of mineral technology inc nysemtx ranger energy service rngr to release earnings on monday hacker exploit microsoft software vulnerability to reportedly target government and businesseswhat to know cwm llc take position in millrose property inc nysemrp acuity inc nyseayi share purchased by cwm llc massachusetts financial service co ma decrease position in select water solution inc nysewttr louisiana state employee retirement system ha stock position in apogee enterprise inc nasdaqapog signaturefd llc lower stock holding in interpublic group of company inc the nyseipg envestnet asset management inc grows holding in yeti holding inc nyseyeti protagonist therapeutic inc nasdaqptgx share sold by massachusetts financial service co ma sbi security co ltd lower position in immunitybio inc nasdaqibrx massachusetts financial service co ma sell share of zurn elkay water solution cor nysezws reynolds consumer product inc nasdaqreyn share sold by massachusetts financial service co ma cwm llc acquires share of spdr sp biotech etf nysearcaxbi edgestream partner lp sell share of the traveler company inc nysetrv sbi security co ltd purchase share of fluor corporation nyseflr st bancorp nasdaqstba rating increased to hold at wall street zen sbi security co ltd ha holding in first trust buywrite income etf nasdaqfthi copa holding sa nysecpa position boosted by massachusetts financial service co ma cemex nysecx cut to hold at wall street zen american asset trust nyseaat upgraded by wall street zen to buy rating massachusetts financial service co ma sell share of kiniksa pharmaceutical international plc nasdaqknsa fg annuity life nysefg upgraded by wall street zen to hold rating brookdale senior living nysebkd downgraded by wall street zen to sell westlake corp nysewlk given consensus rating of moderate buy by brokerage upland software nasdaqupld stock rating upgraded by wall street zen envestnet asset management inc raise stock position in yeti holding inc nyseyeti upland software nasdaqupld rating increased to strongbuy at wall street zen dxc technology nysedxc upgraded by wall street zen to buy rating coupang nysecpng lowered to hold rating by wall street zen dxc technology nysedxc upgraded by wall street zen to buy rating ap business summarybrief at am edt rincon broadcasting buy seven local station from imagicomm rep jonathan l jackson buy international business machine corporation nyseibm share grupo aeroportuario del pacifico nysepac upgraded to buy at wall street zen infineon technology otcmktsifnny stock rating upgraded by wall street zen kennedywilson nysekw upgraded to hold at wall street zen neumora therapeutic inc nasdaqnmra receives consensus rating of hold from analyst insteel industry nyseiiin rating increased to strongbuy at wall street zen world share mostly gain after wall street log a rd straight winning week world share mostly gain after wall street log a rd straight winning week telus digital nysetixt rating lowered to hold at wall street zen nuveen churchill direct lending nysencdl upgraded at wall street zen world share mostly gain after wall street log a rd straight winning week orange otcmktsorany lowered to hold rating by wall street zen wall street zen downgrade irsa inversiones y representaciones nyseirs to buy vista energy nysevist upgraded by wall street zen to hold rating nuveen churchill direct lending nysencdl rating increased to hold at wall street zen telus digital nysetixt downgraded by wall street zen to hold orange otcmktsorany lowered to hold rating by wall street zen gold start the week on the upside a investor weigh trump tariff threat sl green realty nyseslg price target raised to karen kwan sell share of aritzia inc tseatz stock nexxen international nasdaqnexn downgraded by wall street zen to buy nexxen international nasdaqnexn downgraded to buy rating by wall street zen david o watson sell share of apellis pharmaceutical inc nasdaqapls stock victory capital management inc purchase share of sally beauty holding inc nysesbh larson financial group llc acquires share of centerpoint energy inc nysecnp victory capital management inc ha million holding in reddit inc nyserddt larson financial group llc buy share of centerpoint energy inc nysecnp victory capital management inc ha million stock position in highwoods property inc nysehiw microsoft hit with sharepoint attack one version still vulnerable apollon wealth management llc sell share of vaneck oil service etf nysearcaoih apollon wealth management llc grows position in netapp inc nasdaqntap biotechne nasdaqtech and neurogene nasdaqngne financial survey hormel food nysehrl versus china marine food group otcmktscmfo head to head survey sentinelone inc nyse share bought by signaturefd llc sigma planning corp ha million holding in jpmorgan chase co nysejpm cohen capital management inc sell share of jpmorgan chase co nysejpm sigma planning corp ha position in invesco rafi u smallmid etf nasdaqprfz apollon wealth management llc ha position in ppl corporation nyseppl signaturefd llc sell share of illumina inc nasdaqilmn visa inc nysev share sold by edgestream partner lp altman plan dc push to democratize ai economic benefit airline stock to follow now july th blue chip stock worth watching july th doj request for voter data from pa other swing state raise concern among election law expert malcolmjamal warner former cosby show star dy at in costa rica drowning there one thing no one know about trump vladimir tenev sell share of robinhood market inc nasdaqhood stock she got the last seat on a flight next to a stranger theyve been married for year insider selling garrett motion inc nysegtx major shareholder sell share of stock top chinese stock to add to your watchlist july th keith a goldan sell share of syndax pharmaceutical inc nasdaqsndx stock insider selling sintana energy inc cvesei director sell share of stock mineral technology inc nysemtx share sold by bailard inc ford motor nysef upgraded by wall street zen to hold rating performance food group company nysepfgc holding cut by massachusetts financial service co ma echostar corporation nasdaqsats receives consensus recommendation of hold from brokerage the tjx company inc nysetjx receives average pt from brokerage zevra therapeutic nasdaqzvra rating lowered to buy at wall street zen coupang nysecpng downgraded to hold rating by wall street zen irsa inversiones y representaciones nyseirs lowered to buy rating by wall street zen internet initiative japan otcmktsiijiy downgraded to hold rating by wall street zen ofg bancorp nyseofg upgraded by wall street zen to hold rating vista energy nysevist rating increased to hold at wall street zen syrian government start evacuating bedouin family from sweida getting antivaxxers to roll up their sleeve akamai technology nasdaqakam downgraded to hold rating by wall street zen victory capital management inc ha million stock position in yum brand inc nyseyum highwoods property inc nysehiw holding raised by victory capital management inc reviewing hormel food nysehrl and china marine food group otcmktscmfo u envoy double down on support for syria government and criticizes israel intervention larson financial group llc acquires share of ct corporation nysects new parliament seating chart show how labor election win shifted the tectonic plate of federal politics data show more than mbta employee made over k in overtime in cullen frost banker inc invests in sea limited sponsored adr nysese cullen frost banker inc ha holding in howmet aerospace inc nysehwm safehold inc nysesafe given consensus rating of hold by brokerage thomson reuters co nysetri stock position raised by cullen frost banker inc ishares core year usd bond etf nasdaqistb share acquired by cullen frost banker inc apollo global management inc nyseapo share acquired by cullen frost banker inc share in fortis nysefts purchased by cullen frost banker inc costar group inc nasdaqcsgp share acquired by cullen frost banker inc verdence capital advisor llc boost stock position in ishares russell midcap value etf nysearcaiws verdence capital advisor llc ha stake in schwab u midcap etf nysearcaschm verdence capital advisor llc ha stake in medpace holding inc nasdaqmedp verdence capital advisor llc sell share of teck resource ltd nyseteck hillenbrand inc nysehi share acquired by verdence capital advisor llc goldman sachs bdc inc nysegsbd share sold by verdence capital advisor llc verdence capital advisor llc buy new stake in ishares sp value etf nysearcaive tetra tech inc nasdaqttek share acquired by verdence capital advisor llc you can walk between the louvre and the guggenheim in this new art district councilor boston mass and cass drug crisis is worsening the city rat infestation verdence capital advisor llc sell share of target corporation nysetgt bailard inc take position in universal insurance holding inc nyseuve bailard inc ha position in agilent technology inc nysea verdence capital advisor llc buy share of axos financial inc nyseax verdence capital advisor llc purchase new holding in dutch bros inc nysebros verdence capital advisor llc ha holding in millerknoll inc nasdaqmlkn bailard inc grows position in umh property inc nyseumh share in mach natural resource lp nysemnr bought by verdence capital advisor llc victory capital management inc boost stock position in dollar tree inc nasdaqdltr victory capital management inc decrease position in charter communication inc nasdaqchtr larson financial group llc increase position in electronic art inc nasdaqea larson financial group llc raise stock position in airbnb inc nasdaqabnb larson financial group llc sell share of rockwell automation inc nyserok oi glass inc nyseoi receives average recommendation of moderate buy from analyst regret thy name is hawley and murkowski and musk sirius xm holding inc nasdaqsiri share sold by cullen frost banker inc wealthfront adviser llc sell share of ishares russell etf nysearcaiwm wealthfront adviser llc purchase share of the chemours company nysecc wealthfront adviser llc ha holding in first trust senior loan etf nasdaqftsl ast spacemobile inc nasdaqasts share acquired by wealthfront adviser llc pvh corp nysepvh share sold by wealthfront adviser llc emerging market internet etf nysearcaemqq share acquired by wealthfront adviser llc wealthfront adviser llc buy share of f inc nasdaqffiv wealthfront adviser llc boost position in invitation home nyseinvh avantis emerging market equity etf nysearcaavem share purchased by wealthfront adviser llc texas republican want to redraw congressional district in special session inside idaho only maximum security prison where bryan kohberger may spend the rest of his life what to expect a harvard take on the trump administration in federal court vances and marco rubios senate successor avoid gop primary drama independent voter in nyc must count is pakistan moving to join the abraham accord wary democratic establishment keeping new york mamdani at arm length how trump gameplanning for the midterm ai chatbots cocreating a matrix of unreality what i won in a staring contest with a lava heron on san cristobal island vances and marco rubios senate successor avoid gop primary drama bailard inc make new investment in pagaya technology ltd nasdaqpgy share in ge vernova inc nysegev acquired by bailard inc bailard inc sell share of philip morris international inc nysepm promising home improvement stock to research july th bailard inc make new investment in angiodynamics inc nasdaqango best ecommerce stock to watch now july th battery technology stock to research july th nanotechnology stock to consider july th value stock to follow today july th syros pharmaceutical inc nasdaqsyrs receives consensus rating of hold from brokerage massachusetts financial service co ma ha holding in vanguard value etf nysearcavtv arcus bioscience inc nysercus share sold by massachusetts financial service co ma massachusetts financial service co ma sell share of anika therapeutic inc nasdaqanik matson suspends electric vehicle shipment over battery fire concern inside pure fuel partnershipsoverretail strategy sleeper dividend stock to buy for massive upside builder mortgage aid contributing to higher home price morgan stanley say notice on the general meeting of shareholder of the auga group ab entity under restructuring on august to approve the company draft restructuring plan gries financial llc ha million position in procter gamble company the nysepg vanguard growth etf nysearcavug set new year high here why vanguard sp etf nysearcavoo is arrow financial corp largest position gries financial llc ha million stake in bristol myers squibb company nysebmy global bond yield come under pressure obama foundation add outgoing ford foundation president to board of director why opendoor stock is skyrocketing today investopedia jpm high crowding lends risk to stock not involved in ai nucleus security steve carter named a security business innovator hydrograph launch compounding partner program to expand grapheneenhanced thermoplastic klp official deny claim made in civil right lawsuit filed against them lighting precinct reject claim made in civil right lawsuit cornerstone bancorp inc quarterly report june deutsche bank aktiengesellschaft raise american express nyseaxp price target to treasure coast financial planning invests in costco wholesale corporation nasdaqcost share in costco wholesale corporation nasdaqcost purchased by tabor asset management lp costco wholesale corporation nasdaqcost holding trimmed by crestwood advisor group llc im able foundation get donation from firstenergy im able foundation get donation from firstenergy im able foundation get donation from firstenergy rosen law firm encourages vera bradley inc investor to inquire about security class action investigation vra ripple xrp trader who made m earlier this year buy again after month but it not xrp hii summer intern research top use case for ai present finding to company leadership samsung galaxy watch why it the best ai wearable wnba allstar player pay move began well before shirt and sign the pulse business new research suggests rising pr ethic problem with organizational gaslighting best samsung smart tv under r in cnbc daily open investor dismiss trump administration beef with the fed sp hit new high ap technology summarybrief at pm edt delta plane from minneapolis avoids near collision with b bomber iran sends a rocket designed to carry satellite into a suborbital test flight iran sends a rocket designed to carry satellite into a suborbital test flight wex tap ai for faster fsa reimbursement ap business summarybrief at pm edt biogen bet big on billion expansion to boost u drug manufacturing retail speculation is back with a vengeance heat pump industry on the rise global market to expand significantly by senior connection receives innovative achievement award for care express weekly analyst rating update for metlife met hazlett burt watson inc boost stock position in schwab u dividend equity etf nysearcaschd hazlett burt watson inc grows holding in salesforce inc nysecrm verus financial partner inc raise stock position in schwab u dividend equity etf nysearcaschd western montana mental health center wmmhc data breach expose personal information murphy law firm investigates legal claim appeal court order new trial for man convicted in etan patz case animal science expertise honored at international conference it end with u insurer file lawsuit say it doesnt have to pay justin baldonis legal fee jonas brother sell portion of music catalog to dad company i dont know that anybody can love your music more silicon lab announces second quarter earnings webcast head to head review ge aerospace nysege v northrop grumman nysenoc van strum towne inc boost holding in mcdonalds corporation nysemcd analyzing storage computer otcmktssoso teradata nysetdc financial review crh crh v it rival gauzy nasdaqgauz versus kla nasdaqklac head to head comparison netflix inc nasdaqnflx share sold by lewis asset management llc replimune group nasdaqrepl see large volume increase should you buy pacific wealth management sell share of abbott laboratory nyseabt walker dunlop nysewd v rocket company nyserkt headtohead survey zegona communication lonzeg hit new year high following analyst upgrade dwave quantum stock is trading higher monday whats going on rogers returning to bargaining table tomorrow usw local member call for real action toward a fair deal southwest airline announces start date for assigned seating cnn democrat try to dissect themselves with autopsy show theyre riddled with disease vista capital partner inc trim stock position in nike inc nysenke van strum towne inc ha million position in pepsico inc nasdaqpep van strum towne inc sell share of ishares core sp etf nysearcaivv integrity investment advisor llc acquires share of ishares core sp etf nysearcaivv unit apartment project proposed for virginia ave in hollywood integrity investment advisor llc buy new stake in ishares core sp etf nysearcaivv astronomer ceo kiss cam controversy sparked over million in prediction market bet on his ouster whataburger saddle up with the new bacon wrangler double a burger a big and bold a the lone star state rosen a ranked and leading law firm encourages petco health and wellness company inc investor to secure counsel before important deadline in security class action woof critical comparison finwise bancorp nasdaqfinw and texas capital bancshares nasdaqtcbi analyzing kuehne nagel international khngy it competitor rd lewis holding inc acquires share of walmart inc nysewmt comparing rayonier advanced material nyseryam and suzano nysesuz comparing osisko mining otcmktsobnnf aris mining nysearmn analyzing elme community nyseelme sun community nysesui marinemax hzo v it competitor headtohead analysis contrasting williamssonoma wsm and it competitor oneascent wealth management llc increase stake in walmart inc nysewmt reviewing ivanhoe mine ivpaf and it competitor contrasting senstar technology nasdaqsnt applied dna science nasdaqapdn ivanhoe mine ivpaf and it competitor financial comparison headtohead analysis suzano nysesuz versus rayonier advanced material nyseryam lewis asset management llc increase stake in cisco system inc nasdaqcsco unifirst nyseunf and vestis nysevsts financial analysis virgin galactic nysespce v safran otcmktssafry critical contrast amazon ups quantum game with new qubit qpu eye commercial future perception that tariff will be pulled back if market react is alive and well say cfrs patterson skip brittenham hollywood lawyer to the star dy at rogers announces pricing of cash tender offer for canadian dollar debt security rep gilbert ray cisneros jr buy northrop grumman corporation nysenoc stock rep gilbert ray cisneros jr sell prologis inc nysepld stock rep robert bresnahan jr sell off share of arista network inc nyseanet rep gilbert ray cisneros jr purchase share of okta inc nasdaqokta aire advisor llc cut stock position in abbvie inc nyseabbv lauterbach financial advisor llc invests in vanguard sp etf nysearcavoo minnesota lawmaker to resign after being convicted of felony burglary research offer link between burn pit smoke and serious brain injury gop megabills final score t in red ink and million kicked off health insurance cbo say gop megabills final score t in red ink and million kicked off health insurance cbo say trump wade into bryan kohberger murder case demand answer on why he did it there are no explanation sentencing hearing set for exkentucky officer convicted in breonna taylor raid treasure coast financial planning purchase share of costco wholesale corporation nasdaqcost ge vernova gev versus it peer critical contrast tabor asset management lp make new investment in costco wholesale corporation nasdaqcost wealth dimension group ltd ha stake in cocacola company the nyseko trump team throw nyc mayor eric adam under the bus the public deserves to know dcfs withholds detail in illinois girl death miller financial service llc sell share of vanguard value etf nysearcavtv maine need to brace for natural disaster opinion maine need to brace for natural disaster opinion maine need to brace for natural disaster opinion appeal court order new trial for man convicted in etan patz case hunter biden go off on george clooney f him and everybody around him activeduty marine will be leaving los angeles pentagon confirms encountering a bear on the trail yieldmax ultra option income strategy etf nysearcaulty see strong trading volume still a buy contrasting cosan csan it rival rollins nyserol versus limbach nasdaqlmb critical comparison southwest airline announces start date for assigned seating rogers returning to bargaining table tomorrow usw local member call for real action toward a fair deal bruker nasdaqbrkr share gap down whats next irobot nasdaqirbt v robogroup tek otcmktsrobof head to head analysis headtohead contrast kuehne nagel international khngy and it peer srh advisor llc cut holding in the charles schwab corporation nyseschw lewis asset management llc sell share of pfizer inc nysepfe limoneira nasdaqlmnr and forafric global nasdaqafri headtohead review analyzing samsonite group smsey it competitor senstar technology nasdaqsnt and applied dna science nasdaqapdn head to head survey malcolmjamal warner cosby show star dy of drowning at okta nasdaqokta stock unloaded rep robert bresnahan jr travis county commissioner margaret gomez to retire after more than year of service resecurity and braly insurance group announce strategic partnership to drive innovation in cybersecurity and insurance risk management pb and npr are generally unbiased independent of government propaganda and provide key benefit to u democracy this analyst just returned from ukraine frontlines here what he saw judge weighs reality of trump ideological deportation policy a activist crackdown trial end judge weighs reality of trump ideological deportation policy a activist crackdown trial end judge weighs reality of trump ideological deportation policy a activist crackdown trial end top metaverse stock to follow today july st judge weighs reality of trump ideological deportation policy a activist crackdown trial end judge weighs reality of trump ideological deportation policy a activist crackdown trial end lowering the voting age is a call to trust young people with democracy likely new owner of paramount pledge to snuff out bias at cbs lawyer tell fcc that cbs news will practice unbiased journalism judge weighs reality of trump ideological deportation policy a activist crackdown trial end no fk you hunter julio rosa wallop hunter biden for whining about illegals not cleaning his toilet the day niagara fall ran dry here what it looked like highpoint advisor group llc boost stock holding in johnson johnson nysejnj highpoint advisor group llc purchase share of uber technology inc nyseuber hotel stock to add to your watchlist july st highpoint advisor group llc ha million stake in bank of america corporation nysebac promising g stock to add to your watchlist july st peterson wealth service purchase share of johnson johnson nysejnj promising digital medium stock to follow today july st best gold stock to watch today july st verizon communication nysevz release quarterly earnings result beat expectation by eps exxon mobil nysexom braskem nysebak headtohead contrast financial stock to keep an eye on july st johnson johnson nysejnj stock holding decreased by peregrine asset adviser inc this pb station in boston sign is going viral after trump defunded them mark herr communication issue statement regarding the sentencing of anne pramaggiore mastercard incorporated nysema share sold by gries financial llc trump say make indian great again year ago he tried to buy and move the cleveland indian hunter biden suggests ambien contributed to joe bidens poor debate performance hunter biden suggests ambien contributed to joe bidens poor debate performance hunter biden suggests ambien contributed to joe bidens poor debate performance hunter biden suggests ambien contributed to joe bidens poor debate performance hunter biden suggests ambien contributed to joe bidens poor debate performance hunter biden suggests ambien contributed to joe bidens poor debate performance visa inc nysev stock position lowered by miller financial service llc arrow financial corp reduces stock holding in the home depot inc nysehd peterson wealth service purchase share of visa inc nysev jw cole advisor inc purchase share of oracle corporation nyseorcl harvard argues against trump admins funding freeze in court hearing future generation global limited asxfgg insider geoff wilson purchase share d boral capital reiterates buy rating for biomea fusion nasdaqbmea pepsi is betting a popular health trend can get more people to drink it soda america golden age is over letter epstein ghost is haunting trump and america enemy could summon more how did being a sport fan get so expensive wealth dimension group ltd boost position in well fargo company nysewfc environmentalist lawsuit to halt alligator alcatraz filed in wrong court florida official say trump official who hold more than one job malcolmjamal warner actor who starred a theo in the cosby show dead at young kim draw another democratic challenger in california hazlett burt watson inc acquires share of caterpillar inc nysecat hazlett burt watson inc ha million stock holding in lockheed martin corporation nyselmt gop house majority shrink after big beautiful victory dhs release video showing attempted robbery of offduty federal officer in nyc update dare bioscience to host august webinar the dare to playtm difference the winch advisory service llc ha million holding in salesforce inc nysecrm clinton email investigation omitted crucial piece of evidence declassified report find sei select international equity etf nasdaqseie see strong trading volume should you buy head to head comparison hydrofarm holding group nasdaqhyfm one bio otcmktsonbi concord medical service nyseccm versus brookdale senior living nysebkd critical contrast'

$ws.Range("A20").Value = 45860
$ws.Range("A20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B20").Value = 6309.6201171875
$ws.Range("C20").Value = 6316.1201171875
$ws.Range("D20").Value = 6281.7099609375
$ws.Range("E20").Value = 6306.60009765625
$ws.Range("F20").Value = 5662040000
$ws.Range("G20").Value = 0.0006375316336257
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = $headline
$ws.Range("J20").Value = $headline
$ws.Range("K20").Value = 0.126
$ws.Range("L20").Value = 0.048
$ws.Range("M20").Value = 0.826
$ws.Range("N20").Value = 0.9998
$ws.Range("O20").Value = 111
$ws.Range("P20").Value = 94
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 297
$ws.Range("S20").Value = 52
$ws.Range("T20").Value = 379
$ws.Range("U20").Value = 146
$ws.Range("V20").Value = 93
$ws.Range("W20").Value = 28
$ws.Range("X20").Value = 155
